$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Market cap" header (B1) to "MarketCap"
$ws.Range("B1").Value = "MarketCap"

# Update the active selection to B1 (matches the saved cursor position)
$ws.Range("B1").Select()
